{"js": "// Update the two-digit x two-digit multiplication answers in the table.\n// Each old expression is unique in the document, so we can safely\n// search-and-replace each one independently.\nconst replacements = [\n  [\"53\u00d724=1272\", \"78\u00d711=858\"],\n  [\"25\u00d773=1825\", \"17\u00d788=1496\"],\n  [\"67\u00d778=5226\", \"47\u00d718=846\"],\n  [\"81\u00d789=7209\", \"14\u00d752=728\"],\n  [\"74\u00d792=6808\", \"74\u00d795=7030\"],\n  [\"84\u00d751=4284\", \"32\u00d799=3168\"],\n  [\"41\u00d789=3649\", \"85\u00d731=2635\"],\n  [\"85\u00d751=4335\", \"82\u00d769=5658\"],\n  [\"20\u00d742=840\", \"25\u00d717=425\"],\n  [\"96\u00d724=2304\", \"75\u00d746=3450\"],\n  [\"63\u00d764=4032\", \"11\u00d737=407\"],\n  [\"93\u00d783=7719\", \"23\u00d777=1771\"],\n  [\"47\u00d788=4136\", \"80\u00d752=4160\"],\n  [\"32\u00d715=480\", \"56\u00d793=5208\"],\n  [\"25\u00d737=925\", \"35\u00d729=1015\"],\n  [\"66\u00d755=3630\", \"43\u00d774=3182\"],\n  [\"19\u00d754=1026\", \"91\u00d751=4641\"],\n  [\"99\u00d788=8712\", \"19\u00d714=266\"],\n  [\"72\u00d790=6480\", \"29\u00d776=2204\"],\n  [\"70\u00d753=3710\", \"58\u00d714=812\"],\n  [\"36\u00d749=1764\", \"44\u00d726=1144\"],\n  [\"24\u00d784=2016\", \"22\u00d762=1364\"],\n  [\"26\u00d759=1534\", \"41\u00d712=492\"],\n  [\"85\u00d781=6885\", \"44\u00d749=2156\"],\n  [\"70\u00d761=4270\", \"33\u00d766=2178\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit x two-digit multiplication answers in the table.\n# Each old expression is unique in the document, so a simple\n# Find/Replace (wdReplaceAll) for each pair is safe and unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"53\u00d724=1272\", \"78\u00d711=858\"),\n    @(\"25\u00d773=1825\", \"17\u00d788=1496\"),\n    @(\"67\u00d778=5226\", \"47\u00d718=846\"),\n    @(\"81\u00d789=7209\", \"14\u00d752=728\"),\n    @(\"74\u00d792=6808\", \"74\u00d795=7030\"),\n    @(\"84\u00d751=4284\", \"32\u00d799=3168\"),\n    @(\"41\u00d789=3649\", \"85\u00d731=2635\"),\n    @(\"85\u00d751=4335\", \"82\u00d769=5658\"),\n    @(\"20\u00d742=840\", \"25\u00d717=425\"),\n    @(\"96\u00d724=2304\", \"75\u00d746=3450\"),\n    @(\"63\u00d764=4032\", \"11\u00d737=407\"),\n    @(\"93\u00d783=7719\", \"23\u00d777=1771\"),\n    @(\"47\u00d788=4136\", \"80\u00d752=4160\"),\n    @(\"32\u00d715=480\", \"56\u00d793=5208\"),\n    @(\"25\u00d737=925\", \"35\u00d729=1015\"),\n    @(\"66\u00d755=3630\", \"43\u00d774=3182\"),\n    @(\"19\u00d754=1026\", \"91\u00d751=4641\"),\n    @(\"99\u00d788=8712\", \"19\u00d714=266\"),\n    @(\"72\u00d790=6480\", \"29\u00d776=2204\"),\n    @(\"70\u00d753=3710\", \"58\u00d714=812\"),\n    @(\"36\u00d749=1764\", \"44\u00d726=1144\"),\n    @(\"24\u00d784=2016\", \"22\u00d762=1364\"),\n    @(\"26\u00d759=1534\", \"41\u00d712=492\"),\n    @(\"85\u00d781=6885\", \"44\u00d749=2156\"),\n    @(\"70\u00d761=4270\", \"33\u00d766=2178\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $result = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $result) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
